# Updated cryptos list on Fri Oct 11 23:12:11 UTC 2024 with GitHub Actions
# Refreshes Price (D) / Volume(1h) (E) figures for existing coins, and
# shifts the tail of the ranking table: row 27 (Binance-PegBSC-USD) drops
# off, rows 28-51 move up one slot with refreshed figures, and a new
# entry (Stellar) is appended as row 51.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) holds plain text (e.g. "62.506.88", European-style
# thousands separators) in the source data. Excel auto-coerces
# numeric-looking text to a real number on assignment (losing formatting /
# precision), so force text storage, then drop back to the default style
# so we don't leave a stray "Text" number format on the cell.
function Set-TextCell($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$rows = @(
    @{ Row=2; B='Bitcoin'; C='https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'; D='62.506.88'; E='  +3.77%  ' },
    @{ Row=3; B='Ethereum'; C='https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'; D='2.407.18'; E='  +1.09%  ' },
    @{ Row=4; B='TetherUSD'; C='https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'; D='1.00'; E='  +0.35%  ' },
    @{ Row=5; B='BNB'; C='https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'; D='572.73'; E='  +1.93%  ' },
    @{ Row=6; B='Solana'; C='https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; D='145.64'; E='  +4.89%  ' },
    @{ Row=7; B='USDC'; C='https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'; D='0.998'; E='  -0.29%  ' },
    @{ Row=8; B='XRP'; C='https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'; D='0.539'; E='  +1.74%  ' },
    @{ Row=9; B='LidoStakedEther'; C='https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'; D='2.434.26'; E='  +2.25%  ' },
    @{ Row=10; B='Dogecoin'; C='https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; D='0.111'; E='  +4.53%  ' },
    @{ Row=11; B='TRON'; C='https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; D='0.160'; E='  +0.65%  ' },
    @{ Row=12; B='Toncoin'; C='https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D='5.23'; E='  +2.37%  ' },
    @{ Row=13; B='Cardano'; C='https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'; D='0.352'; E='  +3.86%  ' },
    @{ Row=14; B='Avalanche'; C='https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D='27.48'; E='  +6.58%  ' },
    @{ Row=15; B='ShibaInu'; C='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D='0.0000176'; E='  +5.64%  ' },
    @{ Row=16; B='WrappedliquidstakedEther2.0'; C='https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; D='2.883.71'; E='  +2.59%  ' },
    @{ Row=17; B='WrappedBTC'; C='https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; D='62.374.60'; E='  +4.11%  ' },
    @{ Row=18; B='WrappedEther'; C='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D='2.431.21'; E='  +2.56%  ' },
    @{ Row=19; B='Uniswap'; C='https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; D='7.91'; E='  -1.62%  ' },
    @{ Row=20; B='Chainlink'; C='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D='10.93'; E='  +3.70%  ' },
    @{ Row=21; B='BitcoinCash'; C='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D='327.27'; E='  +1.47%  ' },
    @{ Row=22; B='Polkadot'; C='https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D='4.13'; E='  +1.68%  ' },
    @{ Row=23; B='SuiNetwork'; C='https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'; D='2.05'; E='  +11.22%  ' },
    @{ Row=24; B='Dai'; C='https://coinranking.com/coin/MoTuySvg7+dai-dai'; D='0.997'; E='  -0.45%  ' },
    @{ Row=25; B='Litecoin'; C='https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; D='65.57'; E='  +2.19%  ' },
    @{ Row=26; B='Bittensor'; C='https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'; D='626.28'; E='  +11.42%  ' },
    @{ Row=27; B='Aptos'; C='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D='8.48'; E='  +4.20%  ' },
    @{ Row=28; B='PEPE'; C='https://coinranking.com/coin/03WI8NQPF+pepe-pepe'; D='0.0₃0980'; E='  +5.20%  ' },
    @{ Row=29; B='WrappedeETH'; C='https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'; D='2.561.12'; E='  +2.67%  ' },
    @{ Row=30; B='InternetComputer(DFINITY)'; C='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D='8.20'; E='  +2.59%  ' },
    @{ Row=31; B='Fetch.AI'; C='https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'; D='1.41'; E='  +7.20%  ' },
    @{ Row=32; B='Kaspa'; C='https://coinranking.com/coin/V8GxkwWow+kaspa-kas'; D='0.137'; E='  +3.32%  ' },
    @{ Row=33; B='PancakeSwap'; C='https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; D='1.84'; E='  +2.55%  ' },
    @{ Row=34; B='BabyDogeCoin'; C='https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'; D='0.0₆0392'; E='  +36.19%  ' },
    @{ Row=35; B='ImmutableX'; C='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D='1.49'; E='  +2.47%  ' },
    @{ Row=36; B='FirstDigitalUSD'; C='https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'; D='0.996'; E='  -0.32%  ' },
    @{ Row=37; B='NEARProtocol'; C='https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; D='4.75'; E='  +3.88%  ' },
    @{ Row=38; B='PolygonEcosystemToken'; C='https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'; D='0.374'; E='  +1.60%  ' },
    @{ Row=39; B='Monero'; C='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D='151.93'; E='  -1.22%  ' },
    @{ Row=40; B='RenderToken'; C='https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'; D='5.38'; E='  +6.06%  ' },
    @{ Row=41; B='EthereumClassic'; C='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D='18.62'; E='  +2.34%  ' },
    @{ Row=42; B='dogwifhat'; C='https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'; D='2.75'; E='  +12.13%  ' },
    @{ Row=43; B='Stacks'; C='https://coinranking.com/coin/mMPrMcB7+stacks-stx'; D='1.74'; E='  +4.76%  ' },
    @{ Row=44; B='USDe'; C='https://coinranking.com/coin/exbfr2U-0+usde-usde'; D='0.999'; E='  -0.03%  ' },
    @{ Row=45; B='WhiteBITCoin'; C='https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'; D='14.74'; E='  +25.87%  ' },
    @{ Row=46; B='Aave'; C='https://coinranking.com/coin/ixgUfzmLR+aave-aave'; D='144.43'; E='  +2.91%  ' },
    @{ Row=47; B='Filecoin'; C='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D='3.59'; E='  +1.85%  ' },
    @{ Row=48; B='InjectiveProtocol'; C='https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'; D='20.53'; E='  +6.84%  ' },
    @{ Row=49; B='Mantle'; C='https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'; D='0.598'; E='  +1.95%  ' },
    @{ Row=50; B='Hedera'; C='https://coinranking.com/coin/jad286TjB+hedera-hbar'; D='0.0515'; E='  +2.58%  ' },
    @{ Row=51; B='Stellar'; C='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D='0.0917'; E='  +2.03%  ' }
)

foreach ($r in $rows) {
    $ws.Range("B$($r.Row)").Value = $r.B
    $ws.Range("C$($r.Row)").Value = $r.C
    Set-TextCell "D$($r.Row)" $r.D
    $ws.Range("E$($r.Row)").Value = $r.E
}
